$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.375839
$ws.Range("H2").Value = 16.127517
$ws.Range("I2").Value = 0.2354568587499626
$ws.Range("J2").Value = 0.2354568587499626
$ws.Range("O2").Value = 0.3099390012751145
$ws.Range("P2").Value = 0.3099390012751145
$ws.Range("Q2").Value = 0.6446347739533334
$ws.Range("R2").Value = 5.80171296558
$ws.Range("S2").Value = 0.07297726364433911
$ws.Range("T2").Value = 0.07297726364433911

# Row 3
$ws.Range("G3").Value = 5.375839
$ws.Range("H3").Value = 16.127517
$ws.Range("I3").Value = 0.2354568587499626
$ws.Range("J3").Value = 0.2354568587499626
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2669800000000001
$ws.Range("N3").Value = 0.8009400000000001
$ws.Range("O3").Value = 0.6900609987248855
$ws.Range("P3").Value = 0.6900609987248854
$ws.Range("Q3").Value = 1.43524149622
$ws.Range("R3").Value = 12.91717346598
$ws.Range("S3").Value = 0.1624795951056235
$ws.Range("T3").Value = 0.1624795951056234

# Row 4
$ws.Range("I4").Value = 0.007131134316291014
$ws.Range("J4").Value = 0.007131134316291014
$ws.Range("O4").Value = 0.3099390012751145
$ws.Range("P4").Value = 0.3099390012751145
$ws.Range("S4").Value = 0.002210216647949934
$ws.Range("T4").Value = 0.002210216647949934

# Row 5
$ws.Range("I5").Value = 0.007131134316291014
$ws.Range("J5").Value = 0.007131134316291014
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2669800000000001
$ws.Range("N5").Value = 0.8009400000000001
$ws.Range("O5").Value = 0.6900609987248855
$ws.Range("P5").Value = 0.6900609987248854
$ws.Range("Q5").Value = 0.04346825970666667
$ws.Range("R5").Value = 0.3912143373600001
$ws.Range("S5").Value = 0.004920917668341081
$ws.Range("T5").Value = 0.00492091766834108

# Row 6
$ws.Range("G6").Value = 9.994147
$ws.Range("H6").Value = 29.982441
$ws.Range("I6").Value = 0.4377345486919088
$ws.Range("J6").Value = 0.4377345486919088
$ws.Range("O6").Value = 0.3099390012751145
$ws.Range("P6").Value = 0.3099390012751145
$ws.Range("Q6").Value = 1.198431480593333
$ws.Range("R6").Value = 10.78588332534
$ws.Range("S6").Value = 0.1356710088451832
$ws.Range("T6").Value = 0.1356710088451832

# Row 7
$ws.Range("G7").Value = 9.994147
$ws.Range("H7").Value = 29.982441
$ws.Range("I7").Value = 0.4377345486919088
$ws.Range("J7").Value = 0.4377345486919088
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2669800000000001
$ws.Range("N7").Value = 0.8009400000000001
$ws.Range("O7").Value = 0.6900609987248855
$ws.Range("P7").Value = 0.6900609987248854
$ws.Range("Q7").Value = 2.668237366060001
$ws.Range("R7").Value = 24.01413629454
$ws.Range("S7").Value = 0.3020635398467256
$ws.Range("T7").Value = 0.3020635398467255

# Row 8
$ws.Range("G8").Value = 0.7761303333333333
$ws.Range("H8").Value = 2.328391
$ws.Range("I8").Value = 0.03399380269149206
$ws.Range("J8").Value = 0.03399380269149207
$ws.Range("O8").Value = 0.3099390012751145
$ws.Range("P8").Value = 0.3099390012751145
$ws.Range("Q8").Value = 0.0930683753711111
$ws.Range("R8").Value = 0.83761537834
$ws.Range("S8").Value = 0.01053600525574435
$ws.Range("T8").Value = 0.01053600525574435

# Row 9
$ws.Range("G9").Value = 0.7761303333333333
$ws.Range("H9").Value = 2.328391
$ws.Range("I9").Value = 0.03399380269149206
$ws.Range("J9").Value = 0.03399380269149207
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2669800000000001
$ws.Range("N9").Value = 0.8009400000000001
$ws.Range("O9").Value = 0.6900609987248855
$ws.Range("P9").Value = 0.6900609987248854
$ws.Range("Q9").Value = 0.2072112763933334
$ws.Range("R9").Value = 1.86490148754
$ws.Range("S9").Value = 0.02345779743574771
$ws.Range("T9").Value = 0.02345779743574771

# Row 10
$ws.Range("G10").Value = 6.522593333333333
$ws.Range("H10").Value = 19.56778
$ws.Range("I10").Value = 0.2856836555503455
$ws.Range("J10").Value = 0.2856836555503455
$ws.Range("O10").Value = 0.3099390012751145
$ws.Range("P10").Value = 0.3099390012751145
$ws.Range("Q10").Value = 0.7821459085777778
$ws.Range("R10").Value = 7.039313177199999
$ws.Range("S10").Value = 0.08854450688189792
$ws.Range("T10").Value = 0.08854450688189792

# Row 11
$ws.Range("G11").Value = 6.522593333333333
$ws.Range("H11").Value = 19.56778
$ws.Range("I11").Value = 0.2856836555503455
$ws.Range("J11").Value = 0.2856836555503455
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2669800000000001
$ws.Range("N11").Value = 0.8009400000000001
$ws.Range("O11").Value = 0.6900609987248855
$ws.Range("P11").Value = 0.6900609987248854
$ws.Range("Q11").Value = 1.741401968133333
$ws.Range("R11").Value = 15.6726177132
$ws.Range("S11").Value = 0.1971391486684476
$ws.Range("T11").Value = 0.1971391486684476
